$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$evidenceLinks = @{
    2 = "EvidencePack/WinterRelease/TCO/WR-001_TCO_20251229_1102_UI_Accounting_AI_Banking.png"
    3 = "EvidencePack/WinterRelease/TCO/WR-002_TCO_20251229_1102_UI_Sales_Tax_AI.png"
    4 = "EvidencePack/WinterRelease/TCO/WR-003_TCO_20251229_1102_UI_Project_Management_AI.png"
    5 = "EvidencePack/WinterRelease/TCO/WR-004_TCO_20251229_1102_UI_Homepage_Dashboard.png; EvidencePack/WinterRelease/TCO/WR-004_TCO_20251229_1103_UI_Finance_AI_Dashboard.png"
    6 = "EvidencePack/WinterRelease/TCO/WR-005_TCO_20251229_1103_UI_Solutions_Specialist_Feed.png"
    7 = "EvidencePack/WinterRelease/TCO/WR-006_TCO_20251229_1105_UI_Customer_Agent_Leads.png"
    8 = "EvidencePack/WinterRelease/TCO/WR-007_TCO_20251229_1105_UI_Intuit_Intelligence.png"
    9 = "EvidencePack/WinterRelease/TCO/WR-008_TCO_20251229_1105_UI_Conversational_BI.png"
    10 = "EvidencePack/WinterRelease/TCO/WR-009_TCO_20251229_1102_UI_KPIs_Customizados.png"
    11 = "EvidencePack/WinterRelease/TCO/WR-010_TCO_20251229_1102_UI_Dashboards.png"
    12 = "EvidencePack/WinterRelease/TCO/WR-011_TCO_20251229_1105_UI_3P_Data_Integrations.png"
    13 = "EvidencePack/WinterRelease/TCO/WR-012_TCO_20251229_1103_UI_Calculated_Fields_Reports.png"
    14 = "EvidencePack/WinterRelease/TCO/WR-013_TCO_20251229_1103_UI_Management_Reports.png"
    15 = "EvidencePack/WinterRelease/TCO/WR-014_TCO_20251229_1105_UI_Benchmarking.png"
    16 = "EvidencePack/WinterRelease/TCO/WR-015_TCO_20251229_1103_UI_Multi_Entity_Reports.png"
    17 = "EvidencePack/WinterRelease/TCO/WR-016_TCO_20251229_1102_UI_Dimension_Assignment.png"
    18 = "EvidencePack/WinterRelease/TCO/WR-017_TCO_20251229_1103_UI_Hierarchical_Dimensions.png"
    19 = "EvidencePack/WinterRelease/TCO/WR-018_TCO_20251229_1103_UI_Dimensions_Workflow.png"
    20 = "EvidencePack/WinterRelease/TCO/WR-019_TCO_20251229_1103_UI_Dimensions_Balance_Sheet.png"
    21 = "EvidencePack/WinterRelease/TCO/WR-020_TCO_20251229_1103_UI_Parallel_Approval.png"
    22 = "EvidencePack/WinterRelease/TCO/WR-021_TCO_20251229_1105_UI_Desktop_Migration.png"
    23 = "EvidencePack/WinterRelease/TCO/WR-022_TCO_20251229_1105_UI_DFY_Migration.png"
    24 = "EvidencePack/WinterRelease/TCO/WR-023_TCO_20251229_1105_UI_Feature_Compatibility.png"
    25 = "EvidencePack/WinterRelease/TCO/WR-024_TCO_20251229_1103_UI_Certified_Payroll.png"
    26 = "EvidencePack/WinterRelease/TCO/WR-025_TCO_20251229_1105_UI_Sales_Order.png"
    27 = "EvidencePack/WinterRelease/TCO/WR-026_TCO_20251229_1103_UI_Multi_Entity_Payroll.png"
    28 = "EvidencePack/WinterRelease/TCO/WR-027_TCO_20251229_1105_UI_Garnishments.png"
    29 = "EvidencePack/WinterRelease/TCO/WR-028_TCO_20251229_1105_UI_QBTime_Assignments.png"
    30 = "EvidencePack/WinterRelease/TCO/WR-029_TCO_20251229_1105_UI_Enhanced_Amendments.png"
}

foreach ($row in $evidenceLinks.Keys) {
    $ws.Cells.Item($row, 14).Value = $evidenceLinks[$row]
}

Write-Host "Done populating Evidence_links for" $evidenceLinks.Count "rows"
